# Insert a new data row at row 254 (pushing existing rows 254:365 down to
# 255:366) and populate it with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(254).Insert()

$ws.Cells.Item(254, 1).Value  = 3
$ws.Cells.Item(254, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(254, 3).Value  = "Coquimbo"
$ws.Cells.Item(254, 4).Value  = 44452
$ws.Cells.Item(254, 5).Value  = 5
$ws.Cells.Item(254, 6).Value  = 100112024
$ws.Cells.Item(254, 7).Value  = "Choclo"
$ws.Cells.Item(254, 8).Value  = "Dulce o Americano"
$ws.Cells.Item(254, 9).Value  = "Primera"
$ws.Cells.Item(254, 10).Value = 38
$ws.Cells.Item(254, 11).Value = 32000
$ws.Cells.Item(254, 12).Value = 32000
$ws.Cells.Item(254, 13).Value = 32000
$ws.Cells.Item(254, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(254, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(254, 16).Value = 457
$ws.Cells.Item(254, 17).Value = 70
$ws.Cells.Item(254, 18).Value = "Hortaliza"
